# Applies the "Automation HUB" data refresh: inserts two new expense rows
# (IT Support / Professional Services) into the report and reorders the
# pre-existing rows so that the row that used to be #4 now comes right
# after the newly inserted rows, followed by the row that used to be #3.
#
# Columns: A=ID, B=Description, C=Amount, D=Tax, E=Total, F=Currency, G=Date
# Columns A, C, D, E, G hold numeric/date-looking text that must stay as
# literal text (not get auto-converted to numbers/dates), so their
# NumberFormat is forced to "@" (Text) right before the value is assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (1-based) whose content looks numeric/date-like and must be
# written as literal text rather than being auto-converted by Excel.
$textCols = @(1, 3, 4, 5, 7)

function Set-RowValues($rowIndex, $values) {
    for ($c = 1; $c -le $values.Length; $c++) {
        $cell = $ws.Cells.Item($rowIndex, $c)
        if ($textCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $values[$c - 1]
    }
}

# Make room for the two new rows by inserting two blank rows at row 3;
# this pushes the former rows 3-10 down to rows 5-12.
$ws.Rows("3:4").Insert()

# New row 3: IT Support expense dated 2017-03-19
Set-RowValues 3 @("431391", "IT Support", "4714", "942.8", "5656.8", "CAD", "2017-03-19")

# New row 4: Professional Services expense dated 2017-03-14
Set-RowValues 4 @("219659", "Professional Services", "165237", "33047.4", "198284", "EUR", "2017-03-14")

# The former row 4 (814787 / IT Support / 2017-04-25) now lands in row 5,
# and the former row 3 (109291 / IT Support / 2017-07-05) lands in row 6 -
# swap their contents back so the final order matches the source data.
Set-RowValues 5 @("814787", "IT Support", "99526", "19905.2", "119431", "USD", "2017-04-25")
Set-RowValues 6 @("109291", "IT Support", "65262", "13052.4", "78314.4", "CAD", "2017-07-05")
